$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 1158.538592272451

# Row 3
$ws.Range("B3").Value = 0.6606524410359556
$ws.Range("C3").Value = 3286.919754855326
$ws.Range("D3").Value = 401567.231247708
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 405987.8485715312

# Row 4
$ws.Range("B4").Value = 0.6606524410359556
$ws.Range("C4").Value = 10.34677158129881
$ws.Range("D4").Value = 261.3203778131603
$ws.Range("E4").Value = 1133.036916526867
$ws.Range("G4").Value = 1405.364718362363

# Row 5
$ws.Range("B5").Value = 0.6606524410359556
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 22.3905356188092
$ws.Range("E5").Value = 2195978.878461985
$ws.Range("G5").Value = 2196003.585428127
